# Added filtering options for the Component Analysis
# Removes the trailing forecast-horizon error values (columns G:K or a
# subset thereof) from each data row so that every row only keeps the
# quarters-ahead columns that are actually available for that vintage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row number -> first column letter to clear through column K
# (the table's last populated column for rows 2-43; row 44 only went to J).
$firstRemoved = @{
    2  = 'G'
    3  = 'I'
    4  = 'G'
    5  = 'I'
    6  = 'G'
    7  = 'I'
    8  = 'G'
    9  = 'I'
    10 = 'G'
    11 = 'I'
    12 = 'G'
    13 = 'I'
    14 = 'G'
    15 = 'I'
    16 = 'G'
    17 = 'I'
    18 = 'K'
    19 = 'I'
    20 = 'K'
    21 = 'I'
    22 = 'K'
    23 = 'J'
    24 = 'I'
    26 = 'K'
    27 = 'J'
    28 = 'I'
    30 = 'K'
    31 = 'J'
    32 = 'I'
    34 = 'K'
    35 = 'J'
    36 = 'I'
    38 = 'K'
    39 = 'J'
    40 = 'I'
    42 = 'K'
    43 = 'J'
    44 = 'I'
}

foreach ($row in $firstRemoved.Keys) {
    $startCol = $firstRemoved[$row]
    $lastCol = if ($row -eq 44) { 'J' } else { 'K' }
    $range = $ws.Range("$startCol$row`:$lastCol$row")
    $range.ClearContents()
}
